{"js": "// Replace each three-digit-by-one-digit multiplication answer with its\n// updated counterpart. The mapping below reflects the exact old -> new\n// text pairs from the target diff (applied in document order; each old\n// value is unique so plain text search/replace is unambiguous).\nconst replacements = [\n  [\"974\u00d73=2922\", \"625\u00d74=2500\"],\n  [\"136\u00d79=1224\", \"269\u00d79=2421\"],\n  [\"407\u00d73=1221\", \"682\u00d74=2728\"],\n  [\"327\u00d78=2616\", \"113\u00d74=452\"],\n  [\"340\u00d75=1700\", \"703\u00d76=4218\"],\n  [\"742\u00d78=5936\", \"517\u00d79=4653\"],\n  [\"140\u00d77=980\", \"230\u00d78=1840\"],\n  [\"834\u00d75=4170\", \"666\u00d77=4662\"],\n  [\"646\u00d76=3876\", \"468\u00d77=3276\"],\n  [\"537\u00d77=3759\", \"538\u00d74=2152\"],\n  [\"252\u00d77=1764\", \"821\u00d77=5747\"],\n  [\"505\u00d78=4040\", \"894\u00d78=7152\"],\n  [\"305\u00d79=2745\", \"505\u00d75=2525\"],\n  [\"746\u00d79=6714\", \"564\u00d79=5076\"],\n  [\"644\u00d78=5152\", \"280\u00d79=2520\"],\n  [\"791\u00d77=5537\", \"115\u00d72=230\"],\n  [\"221\u00d74=884\", \"205\u00d73=615\"],\n  [\"988\u00d79=8892\", \"973\u00d78=7784\"],\n  [\"374\u00d73=1122\", \"842\u00d75=4210\"],\n  [\"668\u00d79=6012\", \"622\u00d76=3732\"],\n  [\"730\u00d75=3650\", \"812\u00d73=2436\"],\n  [\"334\u00d74=1336\", \"259\u00d72=518\"],\n  [\"166\u00d77=1162\", \"133\u00d75=665\"],\n  [\"277\u00d73=831\", \"482\u00d77=3374\"],\n  [\"140\u00d78=1120\", \"435\u00d75=2175\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication answer with its\n# updated counterpart. The pairs below reflect the exact old -> new text\n# values from the target diff (every old value is unique in the document,\n# so a case-sensitive Find/Replace-All is unambiguous for each one).\n$replacements = @(\n    @(\"974\u00d73=2922\", \"625\u00d74=2500\"),\n    @(\"136\u00d79=1224\", \"269\u00d79=2421\"),\n    @(\"407\u00d73=1221\", \"682\u00d74=2728\"),\n    @(\"327\u00d78=2616\", \"113\u00d74=452\"),\n    @(\"340\u00d75=1700\", \"703\u00d76=4218\"),\n    @(\"742\u00d78=5936\", \"517\u00d79=4653\"),\n    @(\"140\u00d77=980\", \"230\u00d78=1840\"),\n    @(\"834\u00d75=4170\", \"666\u00d77=4662\"),\n    @(\"646\u00d76=3876\", \"468\u00d77=3276\"),\n    @(\"537\u00d77=3759\", \"538\u00d74=2152\"),\n    @(\"252\u00d77=1764\", \"821\u00d77=5747\"),\n    @(\"505\u00d78=4040\", \"894\u00d78=7152\"),\n    @(\"305\u00d79=2745\", \"505\u00d75=2525\"),\n    @(\"746\u00d79=6714\", \"564\u00d79=5076\"),\n    @(\"644\u00d78=5152\", \"280\u00d79=2520\"),\n    @(\"791\u00d77=5537\", \"115\u00d72=230\"),\n    @(\"221\u00d74=884\", \"205\u00d73=615\"),\n    @(\"988\u00d79=8892\", \"973\u00d78=7784\"),\n    @(\"374\u00d73=1122\", \"842\u00d75=4210\"),\n    @(\"668\u00d79=6012\", \"622\u00d76=3732\"),\n    @(\"730\u00d75=3650\", \"812\u00d73=2436\"),\n    @(\"334\u00d74=1336\", \"259\u00d72=518\"),\n    @(\"166\u00d77=1162\", \"133\u00d75=665\"),\n    @(\"277\u00d73=831\", \"482\u00d77=3374\"),\n    @(\"140\u00d78=1120\", \"435\u00d75=2175\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $true, $newText, 2)\n}\n"}
